# Updates cryptos list values (Price / Volume(1h) columns, and the
# coin name/link/price/volume for rows 20-23 which were re-sorted)
# to match the latest scrape, per commit
# "Updated cryptos list on Mon Jun 26 22:44:10 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = '30.233.01'
$ws.Range("E2").Value2 = '  -0.92%  '
$ws.Range("D3").Value2 = '1.857.88'
$ws.Range("E3").Value2 = '  -2.31%  '
$ws.Range("D4").Formula = "'0.9994"
$ws.Range("E4").Value2 = '  -0.05%  '
$ws.Range("D5").Formula = "'235.74"
$ws.Range("E5").Value2 = '  -1.38%  '
$ws.Range("D6").Formula = "'0.9994"
$ws.Range("E6").Value2 = '  -0.01%  '
$ws.Range("D7").Formula = "'0.4780"
$ws.Range("E7").Value2 = '  -2.76%  '
$ws.Range("D8").Formula = "'0.2804"
$ws.Range("E8").Value2 = '  -4.40%  '
$ws.Range("D9").Formula = "'0.06438"
$ws.Range("E9").Value2 = '  -3.93%  '
$ws.Range("D10").Value2 = '1.859.42'
$ws.Range("E10").Value2 = '  -2.25%  '
$ws.Range("D11").Formula = "'0.07390"
$ws.Range("E11").Value2 = '  +0.72%  '
$ws.Range("D12").Formula = "'16.22"
$ws.Range("E12").Value2 = '  -5.07%  '
$ws.Range("D13").Formula = "'5.095"
$ws.Range("E13").Value2 = '  -1.59%  '
$ws.Range("D14").Formula = "'87.08"
$ws.Range("E14").Value2 = '  -1.13%  '
$ws.Range("D15").Formula = "'0.6444"
$ws.Range("E15").Value2 = '  -3.72%  '
$ws.Range("D16").Value2 = '30.177.10'
$ws.Range("E16").Value2 = '  -1.03%  '
$ws.Range("E17").Value2 = '  +0.10%  '
$ws.Range("D18").Formula = "'13.15"
$ws.Range("E18").Value2 = '  -2.56%  '
$ws.Range("D19").Formula = "'0.000007557"
$ws.Range("E19").Value2 = '  -4.30%  '
$ws.Range("B20").Value2 = 'BitcoinCash'
$ws.Range("C20").Value2 = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Formula = "'224.10"
$ws.Range("E20").Value2 = '  +14.13%  '
$ws.Range("B21").Value2 = 'BinanceUSD'
$ws.Range("C21").Value2 = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D21").Formula = "'0.9997"
$ws.Range("E21").Value2 = '  +0.01%  '
$ws.Range("B22").Value2 = 'Uniswap'
$ws.Range("C22").Value2 = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Formula = "'5.272"
$ws.Range("E22").Value2 = '  -3.26%  '
$ws.Range("B23").Value2 = 'BitDAO'
$ws.Range("C23").Value2 = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Range("D23").Formula = "'0.4017"
$ws.Range("E23").Value2 = '  -5.71%  '
$ws.Range("D24").Formula = "'6.078"
$ws.Range("E24").Value2 = '  -0.93%  '
$ws.Range("D25").Formula = "'9.203"
$ws.Range("E25").Value2 = '  -3.38%  '
$ws.Range("D26").Formula = "'163.66"
$ws.Range("E26").Value2 = '  +0.54%  '
$ws.Range("D27").Formula = "'18.49"
$ws.Range("E27").Value2 = '  +0.50%  '
$ws.Range("D28").Formula = "'1.924"
$ws.Range("E28").Value2 = '  -1.07%  '
$ws.Range("D29").Formula = "'1.434"
$ws.Range("E29").Value2 = '  -3.21%  '
$ws.Range("D30").Formula = "'0.09184"
$ws.Range("E30").Value2 = '  +0.04%  '
$ws.Range("D31").Formula = "'4.233"
$ws.Range("E31").Value2 = '  -2.81%  '
$ws.Range("D32").Formula = "'3.944"
$ws.Range("E32").Value2 = '  -3.74%  '
$ws.Range("D33").Formula = "'0.04974"
$ws.Range("E33").Value2 = '  -3.84%  '
$ws.Range("D34").Formula = "'1.146"
$ws.Range("E34").Value2 = '  +3.31%  '
$ws.Range("D35").Formula = "'0.7229"
$ws.Range("E35").Value2 = '  -3.00%  '
$ws.Range("D36").Formula = "'2.687"
$ws.Range("E36").Value2 = '  -1.17%  '
$ws.Range("D37").Formula = "'0.01824"
$ws.Range("E37").Value2 = '  +0.33%  '
$ws.Range("E38").Value2 = '  -3.33%  '
$ws.Range("D39").Formula = "'0.8992"
$ws.Range("E39").Value2 = '  -2.98%  '
$ws.Range("D40").Formula = "'2.034"
$ws.Range("E40").Value2 = '  -1.64%  '
$ws.Range("E41").Value2 = '  -1.34%  '
$ws.Range("D42").Formula = "'105.70"
$ws.Range("E42").Value2 = '  -1.26%  '
$ws.Range("D43").Formula = "'0.4246"
$ws.Range("E43").Value2 = '  -3.64%  '
$ws.Range("D44").Formula = "'0.9999"
$ws.Range("E44").Value2 = '  +0.42%  '
$ws.Range("D45").Formula = "'0.1303"
$ws.Range("E45").Value2 = '  -5.12%  '
$ws.Range("D46").Formula = "'7.269"
$ws.Range("E46").Value2 = '  -4.52%  '
$ws.Range("D47").Formula = "'63.77"
$ws.Range("E47").Value2 = '  -8.33%  '
$ws.Range("D48").Formula = "'1.494"
$ws.Range("E48").Value2 = '  +6.00%  '
$ws.Range("D49").Formula = "'8.670"
$ws.Range("D50").Formula = "'33.63"
$ws.Range("E50").Value2 = '  -4.20%  '
$ws.Range("E51").Value2 = '  -3.55%  '

# The quote-prefixed .Formula assignments above stamp the cell with an
# internal quotePrefix style, distinct from this sheet's default
# (style-less) data cells. Re-apply the original data-cell style across
# the whole data range so no stray per-cell style survives the save.
$ws.Range("B2:E51").Style = $ws.Range("B3").Style
